$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook (after the last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "myInfoDate"

# Populate the data
$newSheet.Range("A1").Value = "date"
$newSheet.Range("B1").Value = 19900409

# Make the new sheet the active sheet/tab
$newSheet.Activate()
$newSheet.Range("A2").Select() | Out-Null
